$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "33069750"
$ws.Range("D16").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E16").Value = "2101"
$ws.Range("F16").Value = 12114

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002372093"
$ws.Range("D17").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E17").Value = "2101"
$ws.Range("F17").Value = 12114

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "33069750"
$ws.Range("D18").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E18").Value = "2102"
$ws.Range("F18").Value = 36341

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1002372093"
$ws.Range("D19").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E19").Value = "2102"
$ws.Range("F19").Value = 36341

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "33069750"
$ws.Range("D20").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E20").Value = "2103"
$ws.Range("F20").Value = 36341

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1002372093"
$ws.Range("D21").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E21").Value = "2103"
$ws.Range("F21").Value = 36341

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "33069750"
$ws.Range("D22").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E22").Value = "2104"
$ws.Range("F22").Value = 36341

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1002372093"
$ws.Range("D23").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E23").Value = "2104"
$ws.Range("F23").Value = 36341

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "33069750"
$ws.Range("D24").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E24").Value = "2105"
$ws.Range("F24").Value = 36341

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1002372093"
$ws.Range("D25").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E25").Value = "2105"
$ws.Range("F25").Value = 36341

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "33069750"
$ws.Range("D26").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E26").Value = "2106"
$ws.Range("F26").Value = 36341

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1002372093"
$ws.Range("D27").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E27").Value = "2106"
$ws.Range("F27").Value = 36341

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "33069750"
$ws.Range("D28").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E28").Value = "2107"
$ws.Range("F28").Value = 36341

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1002372093"
$ws.Range("D29").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E29").Value = "2107"
$ws.Range("F29").Value = 36341

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "33069750"
$ws.Range("D30").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E30").Value = "2108"
$ws.Range("F30").Value = 36341

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1002372093"
$ws.Range("D31").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E31").Value = "2108"
$ws.Range("F31").Value = 36341

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "33069750"
$ws.Range("D32").Value = "KARINA LUCIA GALVIS PINEDA"
$ws.Range("E32").Value = "2109"
$ws.Range("F32").Value = 29073

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "1002372093"
$ws.Range("D33").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E33").Value = "2109"
$ws.Range("F33").Value = 36341

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1002372093"
$ws.Range("D34").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E34").Value = "2110"
$ws.Range("F34").Value = 36341

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1002372093"
$ws.Range("D35").Value = "LILIANA DEL VALLE MULET COMAS"
$ws.Range("E35").Value = "2111"
$ws.Range("F35").Value = 36341
